# Update quarterly recurrence metrics for row 21 (2025Q3)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C21").Value = 210
$ws.Range("D21").Value = 186
$ws.Range("E21").Value = 24
$ws.Range("F21").Value = 53.29512893982808
